# Fix locker/seat-number bug + refresh row 2 and row 3 student data.
# (commit: "locker and seatNumber bug fixed e.g in one locker only one student will be...")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for the two data rows (row 1 is the header row).
$row2 = @{
    A = "1"
    B = "2025-03-07"
    C = "nikhil"
    D = "kjkhuhuj"
    E = "nghkk"
    F = "8651993323"
    G = "06:00-10:00, 22:00-06:00"
    H = "2"
    I = "1"
    J = "350.00"
    K = "50.00"
    L = "1"
    M = "2025-03-08"
}

$row3 = @{
    A = "2"
    B = "2025-02-05"
    C = "satyam"
    D = "satyam"
    E = "ramkrishnanagar"
    F = "7250585057"
    G = "10:00-14:00, 14:00-18:00"
    H = "2"
    I = "1"
    J = "400.00"
    K = "10.00"
    L = "1"
    M = "2025-03-12"
}

$dataRange = $ws.Range("A2:M3")

# Force every cell in the edited range to Text so values like "1", "350.00"
# and date-shaped strings ("2025-03-07") are stored as literal text instead
# of being auto-coerced into numbers/dates by Excel.
$dataRange.NumberFormat = "@"

foreach ($col in $row2.Keys) {
    $ws.Range($col + "2").Value = $row2[$col]
}

foreach ($col in $row3.Keys) {
    $ws.Range($col + "3").Value = $row3[$col]
}

# Drop the temporary Text formatting again so the cells end up back on the
# workbook's default (General) style, matching the original formatting.
$dataRange.ClearFormats()
